# The deck's slide-master theme ("Integral" / "Red Violet" color scheme,
# stored in ppt/theme/theme1.xml) is switched to the stock PowerPoint
# "Office Theme" color scheme (dk2/lt2/accent1-6/hlink/folHlink), matching
# the Office theme already present elsewhere in the package (notes-master
# theme). This mirrors picking a new color theme for the presentation in
# the PowerPoint UI (Design > Variants > Colors), which rewrites the
# slide master's theme color scheme in place.

$p = $ppt.ActivePresentation
$theme = $p.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# Index mapping for ThemeColorScheme (1-based):
#  1 dk1  2 lt1  3 dk2  4 lt2  5 accent1  6 accent2  7 accent3
#  8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink
# RGB() isn't available in this host, so the PowerPoint "long" BGR-packed
# color (R + G*256 + B*65536) is spelled out per value below.
$colors.Colors(1).RGB  = 0        # dk1      000000
$colors.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$colors.Colors(3).RGB  = 6968388  # dk2      44546A
$colors.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$colors.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$colors.Colors(6).RGB  = 3243501  # accent2  ED7D31
$colors.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$colors.Colors(8).RGB  = 49407    # accent4  FFC000
$colors.Colors(9).RGB  = 12874308 # accent5  4472C4
$colors.Colors(10).RGB = 4697456  # accent6  70AD47
$colors.Colors(11).RGB = 12673797 # hlink    0563C1
$colors.Colors(12).RGB = 7491477  # folHlink 954F72
